# "changed in test script"
# Adds a second test case (Test02 - "Test Gmail Login") to the keyword
# driven framework workbook:
#  - TestCases sheet gets a new row for Test02
#  - TestSteps sheet: appends the eight new TestSteps rows for Test02,
#    a trailing "End" marker row, and normalizes the casing of two
#    existing keywords (VerifyPage -> verifyPage, EnterText -> enterText).

$wb = $excel.ActiveWorkbook
$wsCases = $wb.Worksheets.Item("TestCases")
$wsSteps = $wb.Worksheets.Item("TestSteps")

# --- TestCases sheet: add Test02 row ---------------------------------
$wsCases.Range("A3").Value = "Test02"
$wsCases.Range("B3").Value = "Test Gmail Login"
$wsCases.Range("C3").Value = "Yes"

# --- TestSteps sheet: append the "End" marker row ----------------------
$wsSteps.Range("A14").Value = "End"

# --- TestSteps sheet: append Test02 steps, column by column ------------
$wsSteps.Range("A6").Value = "Test02"
$wsSteps.Range("A7").Value = "Test02"
$wsSteps.Range("A8").Value = "Test02"
$wsSteps.Range("A9").Value = "Test02"
$wsSteps.Range("A10").Value = "Test02"
$wsSteps.Range("A11").Value = "Test02"
$wsSteps.Range("A12").Value = "Test02"
$wsSteps.Range("A13").Value = "Test02"

$wsSteps.Range("B6").Value = "TS_001"
$wsSteps.Range("B7").Value = "TS_002"
$wsSteps.Range("B8").Value = "TS_003"
$wsSteps.Range("B9").Value = "TS_004"
$wsSteps.Range("B10").Value = "TS_005"
$wsSteps.Range("B11").Value = "TS_006"
$wsSteps.Range("B12").Value = "TS_007"
$wsSteps.Range("B13").Value = "TS_008"

$wsSteps.Range("C6").Value = "Launch http://www.google.com"
$wsSteps.Range("C7").Value = "Verify that google home page should display"
$wsSteps.Range("C8").Value = "Click on gmail link"
$wsSteps.Range("C9").Value = "verify the gmail login page should display"
$wsSteps.Range("C10").Value = "Enter text in username field"
$wsSteps.Range("C11").Value = "Enter text in password field"
$wsSteps.Range("C12").Value = "Click on Login Button"
$wsSteps.Range("C13").Value = "Verify that user should logged in"

$wsSteps.Range("D6").Value = "launchBrowser"
$wsSteps.Range("D7").Value = "verifyPage"
$wsSteps.Range("D8").Value = "click"
$wsSteps.Range("D9").Value = "verifyPage"
$wsSteps.Range("D10").Value = "enterText"
$wsSteps.Range("D11").Value = "enterText"
$wsSteps.Range("D12").Value = "click"
$wsSteps.Range("D13").Value = "verifyPage"

# --- TestSteps sheet: normalize existing keyword casing -----------------
$wsSteps.Range("D3").Value = "verifyPage"
$wsSteps.Range("D4").Value = "enterText"

# --- Restore the last-active-cell selections seen in the authored file -
$wsCases.Range("A12").Select()
$wsSteps.Range("A15").Select()
